$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Student Summary")

# --- Insert two new rows (11 & 12) of information: Course Code / Max marks ---
# Copy formatting from the existing row 10 (style index 2) down into the new rows
$ws1.Range("A10:C10").Copy()
$ws1.Range("A11:C12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("B11").Value = "Course Code:"
$ws1.Range("C11").Value = "DSPC607"
$ws1.Range("B12").Value = "Max marks"
$ws1.Range("C12").Value = 15

# --- Update the summary statistics table (rows 13-20) ---
$ws1.Range("A17").Value = "Average Marks"
$ws1.Range("A18").Value = "Less Than 40%"
$ws1.Range("A19").Value = "Between 40 % - 75 %"
$ws1.Range("A20").Value = "More than 75%"

$ws1.Range("B15").Value = 48
$ws1.Range("B16").Value = 3
$ws1.Range("B17").Value = 11.71
$ws1.Range("B18").Value = 3
$ws1.Range("B19").Value = 11
$ws1.Range("B20").Value = 28
